# Applies the diff to interaction_scores.xlsx:
#  - Rows 31-56 get updated productId/score/timestamp values (re-ordering /
#    rescoring of existing interaction records).
#  - Rows 57-61 are newly appended interaction records.
#  - Sheet dimension and ignoredErrors range grow from A1:D56 to A1:D61
#    automatically once the new rows are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 31
$lastRow = 61
$rowCount = $lastRow - $firstRow + 1
$data = New-Object 'object[,]' $rowCount,4

$data[0,0] = "6738b019504ed0629a25b8b1"   # A31
$data[0,1] = "67447fbced2b056beb0f8e01"   # B31
$data[0,2] = 0.7   # C31
$data[0,3] = "2025-03-12T11:31:56.846Z"   # D31
$data[1,0] = "6738b019504ed0629a25b8b1"   # A32
$data[1,1] = "676137906c06138b1419f8a5"   # B32
$data[1,2] = 0.7   # C32
$data[1,3] = "2025-03-12T11:31:56.846Z"   # D32
$data[2,0] = "6738b019504ed0629a25b8b1"   # A33
$data[2,1] = "67f08d651841d535b6af6f57"   # B33
$data[2,2] = 0.7   # C33
$data[2,3] = "2025-03-12T11:31:56.846Z"   # D33
$data[3,0] = "6738b019504ed0629a25b8b1"   # A34
$data[3,1] = "6728ebbb071b8fcf4f501e11"   # B34
$data[3,2] = 0.7   # C34
$data[3,3] = "2025-03-12T11:31:56.846Z"   # D34
$data[4,0] = "6738b019504ed0629a25b8b1"   # A35
$data[4,1] = "6728ea62071b8fcf4f501e02"   # B35
$data[4,2] = 0.7   # C35
$data[4,3] = "2025-03-12T11:31:56.846Z"   # D35
$data[5,0] = "6738b019504ed0629a25b8b1"   # A36
$data[5,1] = "6728e9ab071b8fcf4f501df6"   # B36
$data[5,2] = 0.7   # C36
$data[5,3] = "2025-03-12T11:31:56.846Z"   # D36
$data[6,0] = "6738b019504ed0629a25b8b1"   # A37
$data[6,1] = "67250625bb931ab886fc69db"   # B37
$data[6,2] = 0.7   # C37
$data[6,3] = "2025-03-12T11:31:56.846Z"   # D37
$data[7,0] = "6738b019504ed0629a25b8b1"   # A38
$data[7,1] = "6728e9e6071b8fcf4f501dfc"   # B38
$data[7,2] = 0.7   # C38
$data[7,3] = "2025-03-12T11:31:56.846Z"   # D38
$data[8,0] = "6738b019504ed0629a25b8b1"   # A39
$data[8,1] = "6728ea18071b8fcf4f501dff"   # B39
$data[8,2] = 0.7   # C39
$data[8,3] = "2025-03-12T11:31:56.846Z"   # D39
$data[9,0] = "6738b019504ed0629a25b8b1"   # A40
$data[9,1] = "6728eac6071b8fcf4f501e05"   # B40
$data[9,2] = 0.7   # C40
$data[9,3] = "2025-03-12T11:31:56.846Z"   # D40
$data[10,0] = "6738b019504ed0629a25b8b1"   # A41
$data[10,1] = "67f095081841d535b6af6fae"   # B41
$data[10,2] = 0.2   # C41
$data[10,3] = "2025-06-01T04:37:38.841Z"   # D41
$data[11,0] = "6738b019504ed0629a25b8b1"   # A42
$data[11,1] = "67f093ce1841d535b6af6f93"   # B42
$data[11,2] = 0.1   # C42
$data[11,3] = "2025-06-01T04:38:31.695Z"   # D42
$data[12,0] = "6738b019504ed0629a25b8b1"   # A43
$data[12,1] = "67f0959c1841d535b6af6fb4"   # B43
$data[12,2] = 0.2   # C43
$data[12,3] = "2025-06-01T04:49:13.113Z"   # D43
$data[13,0] = "6738b019504ed0629a25b8b1"   # A44
$data[13,1] = "67f094ec1841d535b6af6fab"   # B44
$data[13,2] = 0.1   # C44
$data[13,3] = "2025-06-01T04:52:37.819Z"   # D44
$data[14,0] = "6738b019504ed0629a25b8b1"   # A45
$data[14,1] = "6728e9cd071b8fcf4f501df9"   # B45
$data[14,2] = 0.1   # C45
$data[14,3] = "2025-06-01T04:53:07.969Z"   # D45
$data[15,0] = "6738b019504ed0629a25b8b1"   # A46
$data[15,1] = "67f095dc1841d535b6af6fba"   # B46
$data[15,2] = 0.1   # C46
$data[15,3] = "2025-06-01T05:24:22.234Z"   # D46
$data[16,0] = "6738b019504ed0629a25b8b1"   # A47
$data[16,1] = "67f094111841d535b6af6f99"   # B47
$data[16,2] = 0.1   # C47
$data[16,3] = "2025-06-01T05:27:27.719Z"   # D47
$data[17,0] = "6738b019504ed0629a25b8b1"   # A48
$data[17,1] = "67f092321841d535b6af6f81"   # B48
$data[17,2] = 0.3   # C48
$data[17,3] = "2025-06-01T05:45:55.487Z"   # D48
$data[18,0] = "6738b019504ed0629a25b8b1"   # A49
$data[18,1] = "6728e8a8071b8fcf4f501df0"   # B49
$data[18,2] = 0.1   # C49
$data[18,3] = "2025-06-01T06:14:02.581Z"   # D49
$data[19,0] = "6738b019504ed0629a25b8b1"   # A50
$data[19,1] = "67f095371841d535b6af6fb1"   # B50
$data[19,2] = 0.3   # C50
$data[19,3] = "2025-06-01T06:15:38.184Z"   # D50
$data[20,0] = "682326702fff19d415752f01"   # A51
$data[20,1] = "6728e9e6071b8fcf4f501dfc"   # B51
$data[20,2] = 0.925   # C51
$data[20,3] = "2025-05-23T04:41:14.096Z"   # D51
$data[21,0] = "682326702fff19d415752f01"   # A52
$data[21,1] = "6728e93a071b8fcf4f501df3"   # B52
$data[21,2] = 0.775   # C52
$data[21,3] = "2025-05-23T09:19:25.598Z"   # D52
$data[22,0] = "682326702fff19d415752f01"   # A53
$data[22,1] = "6728ea18071b8fcf4f501dff"   # B53
$data[22,2] = 0.775   # C53
$data[22,3] = "2025-05-23T09:19:36.004Z"   # D53
$data[23,0] = "682326702fff19d415752f01"   # A54
$data[23,1] = "67f091181841d535b6af6f7b"   # B54
$data[23,2] = 0.775   # C54
$data[23,3] = "2025-05-23T09:19:59.350Z"   # D54
$data[24,0] = "682326702fff19d415752f01"   # A55
$data[24,1] = "6728e9cd071b8fcf4f501df9"   # B55
$data[24,2] = 0.925   # C55
$data[24,3] = "2025-05-23T09:21:52.485Z"   # D55
$data[25,0] = "682326702fff19d415752f01"   # A56
$data[25,1] = "6728f96acb86d3695fa1f4a6"   # B56
$data[25,2] = 0.775   # C56
$data[25,3] = "2025-05-23T09:21:00.865Z"   # D56
$data[26,0] = "682326702fff19d415752f01"   # A57
$data[26,1] = "68067dd1286f80e4174d8736"   # B57
$data[26,2] = 1   # C57
$data[26,3] = "2025-05-23T09:23:22.245Z"   # D57
$data[27,0] = "683b37622eb85e2df9802771"   # A58
$data[27,1] = "6728e93a071b8fcf4f501df3"   # B58
$data[27,2] = 1   # C58
$data[27,3] = "2025-05-31T17:08:32.178Z"   # D58
$data[28,0] = "683b37622eb85e2df9802771"   # A59
$data[28,1] = "6728ea18071b8fcf4f501dff"   # B59
$data[28,2] = 0.85   # C59
$data[28,3] = "2025-05-31T17:10:43.229Z"   # D59
$data[29,0] = "683b37622eb85e2df9802771"   # A60
$data[29,1] = "6728f96acb86d3695fa1f4a6"   # B60
$data[29,2] = 0.15   # C60
$data[29,3] = "2025-06-01T02:53:50.928Z"   # D60
$data[30,0] = "683b37622eb85e2df9802771"   # A61
$data[30,1] = "6728ec07071b8fcf4f501e17"   # B61
$data[30,2] = 0.15   # C61
$data[30,3] = "2025-06-01T07:31:20.015Z"   # D61

$ws.Range("A31:D61").Value = $data
